# Auto-generated cell updates derived from the OOXML diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.33"
$ws.Range("E2").Value = "'0.58%"
$ws.Range("D3").Value = "'44.10"
$ws.Range("E3").Value = "'0.62%"
$ws.Range("D4").Value = "'5.574"
$ws.Range("E4").Value = "'2.02%"
$ws.Range("E5").Value = "'0.24%"
$ws.Range("D6").Value = "'1.976"
$ws.Range("E6").Value = "'5.00%"
$ws.Range("E7").Value = "'0.72%"
$ws.Range("D8").Value = "'0.9533"
$ws.Range("E8").Value = "'1.60%"
$ws.Range("D10").Value = "'0.1159"
$ws.Range("E10").Value = "'-2.34%"
$ws.Range("D11").Value = "'0.1869"
$ws.Range("E11").Value = "'-1.42%"
$ws.Range("D12").Value = "'12.60"
$ws.Range("E12").Value = "'45.79%"
$ws.Range("D13").Value = "'0.09940"
$ws.Range("E13").Value = "'3.93%"
$ws.Range("D14").Value = "'0.04716"
$ws.Range("E14").Value = "'14.87%"
$ws.Range("D15").Value = "'0.1068"
$ws.Range("E15").Value = "'-0.04%"
$ws.Range("D16").Value = "'0.001285"
$ws.Range("E16").Value = "'0.33%"
$ws.Range("D17").Value = "'0.04234"
$ws.Range("E17").Value = "'-3.01%"
$ws.Range("D18").Value = "'0.005941"
$ws.Range("E18").Value = "'0.16%"
$ws.Range("D19").Value = "'3.374"
$ws.Range("E19").Value = "'-5.68%"
$ws.Range("D21").Value = "'0.1410"
$ws.Range("E21").Value = "'4.57%"
$ws.Range("E22").Value = "'0.46%"
$ws.Range("D23").Value = "'0.001257"
$ws.Range("E23").Value = "'1.67%"
$ws.Range("D24").Value = "'0.004366"
$ws.Range("E24").Value = "'1.32%"
$ws.Range("E25").Value = "'-3.56%"
$ws.Range("E26").Value = "'-0.64%"
$ws.Range("D38").Value = "'0.02662"
$ws.Range("E38").Value = "'0.23%"
$ws.Range("D39").Value = "'0.05549"
$ws.Range("E39").Value = "'2.42%"
$ws.Range("D40").Value = "'0.007582"
$ws.Range("E40").Value = "'-0.42%"
$ws.Range("D41").Value = "'0.1410"
$ws.Range("E41").Value = "'1.44%"
$ws.Range("D42").Value = "'0.008089"
$ws.Range("E42").Value = "'-19.13%"
$ws.Range("E43").Value = "'-3.65%"
$ws.Range("D44").Value = "'0.008905"
$ws.Range("E44").Value = "'-10.22%"
$ws.Range("D45").Value = "'0.00007259"
$ws.Range("E45").Value = "'5.50%"
$ws.Range("E46").Value = "'-0.31%"
$ws.Range("D47").Value = "'0.004618"
$ws.Range("E47").Value = "'29.60%"
$ws.Range("E48").Value = "'-0.30%"
$ws.Range("E49").Value = "'-0.31%"
$ws.Range("E50").Value = "'-0.31%"
